$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated pl_mw values for 380 kV case (rows 2-25, columns B,D,E,F,G,J,K,L,M,N,O)
# Row 2
$ws.Range("B2").Value = 1.047363247226826
$ws.Range("D2").Value = 0.1940990992983416
$ws.Range("E2").Value = 0.2285010784184545
$ws.Range("F2").Value = 1.832219175787728
$ws.Range("G2").Value = 0.002494117940094704
$ws.Range("J2").Value = 0.3153605015601499
$ws.Range("K2").Value = 0.3618362577875587
$ws.Range("L2").Value = 0.1379720860992535
$ws.Range("M2").Value = 0.2201592650527857
$ws.Range("N2").Value = 2.562446059108983
$ws.Range("O2").Value = 4.525278713370056

# Row 3
$ws.Range("B3").Value = 1.028272804076835
$ws.Range("D3").Value = 0.1941849182040514
$ws.Range("E3").Value = 0.2296497882683779
$ws.Range("F3").Value = 1.838472239546746
$ws.Range("G3").Value = 0.00249658769447901
$ws.Range("J3").Value = 0.3168055990229846
$ws.Range("K3").Value = 0.3202111586456056
$ws.Range("L3").Value = 0.1319910143878857
$ws.Range("M3").Value = 0.2146039941685096
$ws.Range("N3").Value = 2.584557899823347
$ws.Range("O3").Value = 4.542511771281596

# Row 4
$ws.Range("B4").Value = 1.016973291525574
$ws.Range("D4").Value = 0.1942910189364397
$ws.Range("E4").Value = 0.230401655981928
$ws.Range("F4").Value = 1.843074078533647
$ws.Range("G4").Value = 0.002498186931431273
$ws.Range("J4").Value = 0.3177422708571176
$ws.Range("K4").Value = 0.2946534944191512
$ws.Range("L4").Value = 0.1283638274466057
$ws.Range("M4").Value = 0.2112779821511523
$ws.Range("N4").Value = 2.598828927627177
$ws.Range("O4").Value = 4.555070441406798

# Row 5
$ws.Range("B5").Value = 1.012475373485699
$ws.Range("D5").Value = 0.1943477343420597
$ws.Range("E5").Value = 0.2307197797874903
$ws.Range("F5").Value = 1.845141335098226
$ws.Range("G5").Value = 0.002498859514576418
$ws.Range("J5").Value = 0.3181364141913154
$ws.Range("K5").Value = 0.2842392006988632
$ws.Range("L5").Value = 0.1268972113161269
$ws.Range("M5").Value = 0.2099441089053542
$ws.Range("N5").Value = 2.604819165332905
$ws.Range("O5").Value = 4.560685942259227

# Row 6
$ws.Range("B6").Value = 1.01173496077314
$ws.Range("D6").Value = 0.1943579673934508
$ws.Range("E6").Value = 0.2307733132963827
$ws.Range("F6").Value = 1.84549620333658
$ws.Range("G6").Value = 0.002498972459476462
$ws.Range("J6").Value = 0.3182026136181477
$ws.Range("K6").Value = 0.2825099722945197
$ws.Range("L6").Value = 0.1266543785111836
$ws.Range("M6").Value = 0.2097239232469654
$ws.Range("N6").Value = 2.605824392799129
$ws.Range("O6").Value = 4.561648468411164

# Row 7
$ws.Range("B7").Value = 1.016912198127528
$ws.Range("D7").Value = 0.1942917291809891
$ws.Range("E7").Value = 0.2304058987765298
$ws.Range("F7").Value = 1.843101180668697
$ws.Range("G7").Value = 0.002498195917453064
$ws.Range("J7").Value = 0.3177475360003474
$ws.Range("K7").Value = 0.2945130401595861
$ws.Range("L7").Value = 0.1283440014558224
$ws.Range("M7").Value = 0.2112599058126641
$ws.Range("N7").Value = 2.598909006671805
$ws.Range("O7").Value = 4.555144158044584

# Row 8
$ws.Range("B8").Value = 1.040693599327597
$ws.Range("D8").Value = 0.1941176338501478
$ws.Range("E8").Value = 0.2288875060515272
$ws.Range("F8").Value = 1.834217137000827
$ws.Range("G8").Value = 0.002494952363426184
$ws.Range("J8").Value = 0.3158485388194396
$ws.Range("K8").Value = 0.3474841952808276
$ws.Range("L8").Value = 0.1359004994121307
$ws.Range("M8").Value = 0.2182262606427017
$ws.Range("N8").Value = 2.569926187717686
$ws.Range("O8").Value = 4.530810524410981

# Row 9
$ws.Range("B9").Value = 1.090655531594507
$ws.Range("D9").Value = 0.1941979692905562
$ws.Range("E9").Value = 0.2262782120474434
$ws.Range("F9").Value = 1.822835173832935
$ws.Range("G9").Value = 0.002489245948122587
$ws.Range("J9").Value = 0.3125151916360136
$ws.Range("K9").Value = 0.4513423616591865
$ws.Range("L9").Value = 0.1510729820181069
$ws.Range("M9").Value = 0.2325558483941137
$ws.Range("N9").Value = 2.518594418045847
$ws.Range("O9").Value = 4.498764418169145

# Row 10
$ws.Range("B10").Value = 1.129365175129124
$ws.Range("D10").Value = 0.1945115289514305
$ws.Range("E10").Value = 0.2245841124901569
$ws.Range("F10").Value = 1.818142933660752
$ws.Range("G10").Value = 0.002485448362947259
$ws.Range("J10").Value = 0.3103025940868873
$ws.Range("K10").Value = 0.5276155879779196
$ws.Range("L10").Value = 0.1624312376996784
$ws.Range("M10").Value = 0.2434853567258344
$ws.Range("N10").Value = 2.484227585330357
$ws.Range("O10").Value = 4.484753263684127

# Row 11
$ws.Range("B11").Value = 1.147404941695214
$ws.Range("D11").Value = 0.1947089140672276
$ws.Range("E11").Value = 0.2238615035061651
$ws.Range("F11").Value = 1.816802847384992
$ws.Range("G11").Value = 0.002483805669509031
$ws.Range("J11").Value = 0.3093469997838643
$ws.Range("K11").Value = 0.5623035555405238
$ws.Range("L11").Value = 0.1676432578062759
$ws.Range("M11").Value = 0.2485434656901262
$ws.Range("N11").Value = 2.469318109277642
$ws.Range("O11").Value = 4.480445111492628

# Row 12
$ws.Range("B12").Value = 1.154297514694008
$ws.Range("D12").Value = 0.1947914841178573
$ws.Range("E12").Value = 0.2235947539111853
$ws.Range("F12").Value = 1.8164094097434
$ws.Range("G12").Value = 0.002483195762388349
$ws.Range("J12").Value = 0.3089924377719679
$ws.Range("K12").Value = 0.5754371771199089
$ws.Range("L12").Value = 0.1696232899768404
$ws.Range("M12").Value = 0.2504711097603973
$ws.Range("N12").Value = 2.46377631726255
$ws.Range("O12").Value = 4.479110356870905

# Row 13
$ws.Range("B13").Value = 1.15281035891158
$ws.Range("D13").Value = 0.1947733539379044
$ws.Range("E13").Value = 0.2236518972703623
$ws.Range("F13").Value = 1.816489075758241
$ws.Range("G13").Value = 0.002483326577592793
$ws.Range("J13").Value = 0.309068474695767
$ws.Range("K13").Value = 0.5726087141951268
$ws.Range("L13").Value = 0.1691965741440811
$ws.Range("M13").Value = 0.2500554147838159
$ws.Range("N13").Value = 2.464965211790842
$ws.Range("O13").Value = 4.479384632666438

# Row 14
$ws.Range("B14").Value = 1.147970772667264
$ws.Range("D14").Value = 0.1947155506105389
$ws.Range("E14").Value = 0.2238394199620704
$ws.Range("F14").Value = 1.816768195055744
$ws.Range("G14").Value = 0.002483755248817679
$ws.Range("J14").Value = 0.3093176835627132
$ws.Range("K14").Value = 0.5633841094751233
$ws.Range("L14").Value = 0.1678060296661954
$ws.Range("M14").Value = 0.2487018095220961
$ws.Range("N14").Value = 2.468860096504289
$ws.Range("O14").Value = 4.480329357314901

# Row 15
$ws.Range("B15").Value = 1.145014352003585
$ws.Range("D15").Value = 0.1946811619833753
$ws.Range("E15").Value = 0.2239551792284589
$ws.Range("F15").Value = 1.816954006584993
$ws.Range("G15").Value = 0.002484019402542766
$ws.Range("J15").Value = 0.3094712813067506
$ws.Range("K15").Value = 0.5577334983305775
$ws.Range("L15").Value = 0.1669551045173705
$ws.Range("M15").Value = 0.2478742777973579
$ws.Range("N15").Value = 2.471259382250764
$ws.Range("O15").Value = 4.480946649624883

# Row 16
$ws.Range("B16").Value = 1.128194845826499
$ws.Range("D16").Value = 0.1944997269028832
$ws.Range("E16").Value = 0.22463230158837
$ws.Range("F16").Value = 1.818246479163363
$ws.Range("G16").Value = 0.002485557419528423
$ws.Range("J16").Value = 0.3103660672844661
$ws.Range("K16").Value = 0.5253484089832341
$ws.Range("L16").Value = 0.1620915161097258
$ws.Range("M16").Value = 0.2431565195306504
$ws.Range("N16").Value = 2.485216524008639
$ws.Range("O16").Value = 4.485076344018495

# Row 17
$ws.Range("B17").Value = 1.117986466268405
$ws.Range("D17").Value = 0.1944024099422492
$ws.Range("E17").Value = 0.2250599840169052
$ws.Range("F17").Value = 1.81924269053458
$ws.Range("G17").Value = 0.002486522634949468
$ws.Range("J17").Value = 0.3109280171694575
$ws.Range("K17").Value = 0.5054784115311577
$ws.Range("L17").Value = 0.1591193182569981
$ws.Range("M17").Value = 0.2402843038055948
$ws.Range("N17").Value = 2.493964231918281
$ws.Range("O17").Value = 4.488138554595025

# Row 18
$ws.Range("B18").Value = 1.112155438651968
$ws.Range("D18").Value = 0.1943515915288714
$ws.Range("E18").Value = 0.2253104992136512
$ws.Range("F18").Value = 1.819890474016731
$ws.Range("G18").Value = 0.002487085790743094
$ws.Range("J18").Value = 0.3112560307427077
$ws.Range("K18").Value = 0.4940488963437701
$ws.Range("L18").Value = 0.1574140408273621
$ws.Range("M18").Value = 0.238640405794257
$ws.Range("N18").Value = 2.499063845996414
$ws.Range("O18").Value = 4.490094303151153

# Row 19
$ws.Range("B19").Value = 1.110188138046055
$ws.Range("D19").Value = 0.1943352722437055
$ws.Range("E19").Value = 0.2253960969595088
$ws.Range("F19").Value = 1.820122654025624
$ws.Range("G19").Value = 0.00248727783952168
$ws.Range("J19").Value = 0.3113679147131219
$ws.Range("K19").Value = 0.490178938141554
$ws.Range("L19").Value = 0.1568373974947548
$ws.Range("M19").Value = 0.2380852101498974
$ws.Range("N19").Value = 2.500802195922608
$ws.Range("O19").Value = 4.490789893909692

# Row 20
$ws.Range("B20").Value = 1.119068972061655
$ws.Range("D20").Value = 0.1944122362469329
$ws.Range("E20").Value = 0.225013988480165
$ws.Range("F20").Value = 1.81912890330571
$ws.Range("G20").Value = 0.002486419059561736
$ws.Range("J20").Value = 0.3108677005584148
$ws.Range("K20").Value = 0.5075936960635374
$ws.Range("L20").Value = 0.1594352746585912
$ws.Range("M20").Value = 0.2405892163068799
$ws.Range("N20").Value = 2.493025969010783
$ws.Range("O20").Value = 4.487792455371846

# Row 21
$ws.Range("B21").Value = 1.149390617512381
$ws.Range("D21").Value = 0.1947323168726172
$ws.Range("E21").Value = 0.2237841532733107
$ws.Range("F21").Value = 1.816683118089145
$ws.Range("G21").Value = 0.002483629008470917
$ws.Range("J21").Value = 0.3092442869138265
$ws.Range("K21").Value = 0.5660936578333349
$ws.Range("L21").Value = 0.1682142949045158
$ws.Range("M21").Value = 0.2490990649919311
$ws.Range("N21").Value = 2.467713249361881
$ws.Range("O21").Value = 4.48004382084406

# Row 22
$ws.Range("B22").Value = 1.169564615204251
$ws.Range("D22").Value = 0.1949870937194476
$ws.Range("E22").Value = 0.2230205159191305
$ws.Range("F22").Value = 1.81574917513781
$ws.Range("G22").Value = 0.00248187631700242
$ws.Range("J22").Value = 0.3082258366723964
$ws.Range("K22").Value = 0.6043150446333527
$ws.Range("L22").Value = 0.1739888630590087
$ws.Range("M22").Value = 0.2547320638073671
$ws.Range("N22").Value = 2.451776703163286
$ws.Range("O22").Value = 4.476708562022452

# Row 23
$ws.Range("B23").Value = 1.158764903956808
$ws.Range("D23").Value = 0.1948469591607207
$ws.Range("E23").Value = 0.2234244188273156
$ws.Range("F23").Value = 1.816186905585866
$ws.Range("G23").Value = 0.002482805304347702
$ws.Range("J23").Value = 0.3087655172949697
$ws.Range("K23").Value = 0.5839168583581511
$ws.Range("L23").Value = 0.170903527801201
$ws.Range("M23").Value = 0.2517191501933951
$ws.Range("N23").Value = 2.460226827092338
$ws.Range("O23").Value = 4.478330581590484

# Row 24
$ws.Range("B24").Value = 1.118579453073949
$ws.Range("D24").Value = 0.1944077777913051
$ws.Range("E24").Value = 0.2250347686239698
$ws.Range("F24").Value = 1.819180112725277
$ws.Range("G24").Value = 0.00248646586022698
$ws.Range("J24").Value = 0.3108949543061277
$ws.Range("K24").Value = 0.5066373945328451
$ws.Range("L24").Value = 0.1592924199237444
$ws.Range("M24").Value = 0.2404513423751098
$ws.Range("N24").Value = 2.493449938249795
$ws.Range("O24").Value = 4.487948318630345

# Row 25
$ws.Range("B25").Value = 1.076785989695537
$ws.Range("D25").Value = 0.1941313522831862
$ws.Range("E25").Value = 0.2269448310274962
$ws.Range("F25").Value = 1.825268937525017
$ws.Range("G25").Value = 0.002490720053650364
$ws.Range("J25").Value = 0.3133753120703151
$ws.Range("K25").Value = 0.4232500137970874
$ws.Range("L25").Value = 0.1469310096170631
$ws.Range("M25").Value = 0.2286083616364785
$ws.Range("N25").Value = 2.53189272337162
$ws.Range("O25").Value = 4.505758110286592
